$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.403.77'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.35%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.780.32'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +3.83%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('E5').Value = '  +1.24%  '
$ws.Range('E6').Value = '  -0.09%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5289'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +11.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3765'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +9.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '42.82'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07411'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.34%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.094'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '20.71'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +4.83%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.111'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +4.98%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '1.782.18'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.50%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.989'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.59%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '89.79'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +3.06%  '
$ws.Range('E18').Value = '  +2.27%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06433'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.97%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.08%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.80'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.904'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.444.96'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +1.28%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.65%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.099'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '155.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.377'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +15.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.989.16'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +3.82%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '121.47'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.58%  '
$ws.Range('E31').Value = '  +6.61%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1026'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +11.19%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.594'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.30%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.635'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02261'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05981'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +2.46%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '11.30'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.53%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.920'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.58%  '
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6140'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '8.235'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.434'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('E43').Value = '  +4.69%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.29'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.5797'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +4.20%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.621'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.06%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '121.68'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.40%  '
$ws.Range('E48').Value = '  +3.96%  '
$ws.Range('E49').Value = '  +1.59%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.03'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.70%  '
